{"js": "const replacements = [\n  [\"2024-02-07 Wednesday\", \"2024-02-08 Thursday\"],\n  [\"72\\u00d726=1872\", \"61\\u00d750=3050\"],\n  [\"64\\u00d724=1536\", \"77\\u00d771=5467\"],\n  [\"34\\u00d790=3060\", \"46\\u00d764=2944\"],\n  [\"13\\u00d790=1170\", \"23\\u00d784=1932\"],\n  [\"32\\u00d728=896\", \"30\\u00d763=1890\"],\n  [\"11\\u00d789=979\", \"17\\u00d718=306\"],\n  [\"64\\u00d727=1728\", \"30\\u00d741=1230\"],\n  [\"48\\u00d750=2400\", \"84\\u00d730=2520\"],\n  [\"11\\u00d782=902\", \"62\\u00d739=2418\"],\n  [\"89\\u00d723=2047\", \"55\\u00d759=3245\"],\n  [\"36\\u00d742=1512\", \"80\\u00d722=1760\"],\n  [\"99\\u00d755=5445\", \"59\\u00d728=1652\"],\n  [\"50\\u00d735=1750\", \"69\\u00d724=1656\"],\n  [\"60\\u00d780=4800\", \"70\\u00d780=5600\"],\n  [\"98\\u00d717=1666\", \"17\\u00d795=1615\"],\n  [\"40\\u00d728=1120\", \"18\\u00d778=1404\"],\n  [\"17\\u00d742=714\", \"17\\u00d763=1071\"],\n  [\"11\\u00d751=561\", \"27\\u00d717=459\"],\n  [\"23\\u00d722=506\", \"72\\u00d748=3456\"],\n  [\"83\\u00d723=1909\", \"76\\u00d798=7448\"],\n  [\"30\\u00d727=810\", \"91\\u00d737=3367\"],\n  [\"71\\u00d720=1420\", \"81\\u00d750=4050\"],\n  [\"67\\u00d751=3417\", \"77\\u00d755=4235\"],\n  [\"75\\u00d755=4125\", \"36\\u00d776=2736\"],\n  [\"25\\u00d745=1125\", \"72\\u00d724=1728\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-02-07 Wednesday\", \"2024-02-08 Thursday\"),\n    @(\"72\u00d726=1872\", \"61\u00d750=3050\"),\n    @(\"64\u00d724=1536\", \"77\u00d771=5467\"),\n    @(\"34\u00d790=3060\", \"46\u00d764=2944\"),\n    @(\"13\u00d790=1170\", \"23\u00d784=1932\"),\n    @(\"32\u00d728=896\", \"30\u00d763=1890\"),\n    @(\"11\u00d789=979\", \"17\u00d718=306\"),\n    @(\"64\u00d727=1728\", \"30\u00d741=1230\"),\n    @(\"48\u00d750=2400\", \"84\u00d730=2520\"),\n    @(\"11\u00d782=902\", \"62\u00d739=2418\"),\n    @(\"89\u00d723=2047\", \"55\u00d759=3245\"),\n    @(\"36\u00d742=1512\", \"80\u00d722=1760\"),\n    @(\"99\u00d755=5445\", \"59\u00d728=1652\"),\n    @(\"50\u00d735=1750\", \"69\u00d724=1656\"),\n    @(\"60\u00d780=4800\", \"70\u00d780=5600\"),\n    @(\"98\u00d717=1666\", \"17\u00d795=1615\"),\n    @(\"40\u00d728=1120\", \"18\u00d778=1404\"),\n    @(\"17\u00d742=714\", \"17\u00d763=1071\"),\n    @(\"11\u00d751=561\", \"27\u00d717=459\"),\n    @(\"23\u00d722=506\", \"72\u00d748=3456\"),\n    @(\"83\u00d723=1909\", \"76\u00d798=7448\"),\n    @(\"30\u00d727=810\", \"91\u00d737=3367\"),\n    @(\"71\u00d720=1420\", \"81\u00d750=4050\"),\n    @(\"67\u00d751=3417\", \"77\u00d755=4235\"),\n    @(\"75\u00d755=4125\", \"36\u00d776=2736\"),\n    @(\"25\u00d745=1125\", \"72\u00d724=1728\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute([ref]$old, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$new, [ref]2)\n}\n"}
